$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 8 (pushes "About Page changes" block and everything
# below it down by 5 rows) to make room for the new "News story progress bar" feature.
$ws.Rows("8:12").Insert()

# Fill in the new feature block.
$ws.Range("A8").Value = "News story progress bar"
$ws.Range("B8").Value = "Research"
$ws.Range("C8").Value = 2
$ws.Range("E8").Value = "Determine if needed at all"

$ws.Range("B9").Value = "Codepen design"
$ws.Range("C9").Value = 3

$ws.Range("B10").Value = "Accessibility check"
$ws.Range("C10").Value = 1

$ws.Range("B11").Value = "Test on Mura"
$ws.Range("C11").Value = 2

$ws.Range("B12").Value = "Implement"
$ws.Range("C12").Value = 1

# Grow the worksheet table (ListObject) to cover the newly inserted rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E36"))

# Move the active selection to A2, matching the post-edit workbook state.
$ws.Range("A2").Select() | Out-Null
